# edit.ps1 - applies the tracked change:
#   1) "От  Московского" + " педагогического" (two runs, wrapped in
#      w:proofErr gramStart/gramEnd) -> single run
#      "От  Московского педагогического" (grammar-check markers cleared
#      by the retype).
#   2) "lastNameRU" -> "lastNameR" + "u" (the trailing "U" was retyped as
#      lowercase "u", leaving the text split across two runs).
#
# wdReplaceOne = 2 (Find.Execute Replace: parameter)

$d = $word.ActiveDocument

# --- 1) "От  Московского" / " педагогического" -> merged single run ---
$r1 = $d.Content
$r1.Find.Execute(
    "От  Московского педагогического", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "От  Московского педагогического", 2) | Out-Null

# --- 2) "lastNameRU" -> "lastNameR" run + "u" run ---
$r2 = $d.Content
$r2.Find.Execute(
    "lastNameRU", $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null

$matchEnd = $r2.End

# Replace just the trailing "U" with a lower-case "u" ...
$lastChar = $d.Range($matchEnd - 1, $matchEnd)
$lastChar.Text = "u"

# ... then touch its character formatting so the retyped character keeps
# its own run instead of being re-coalesced with the "lastNameR" run that
# precedes it (mirrors the two distinct <w:r> elements left behind by a
# real Word edit session).
$newChar = $d.Range($matchEnd - 1, $matchEnd)
$newChar.Italic = 0

"done"
